$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.472.88"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.873.60"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").Value = "242.94"
$ws.Range("E5").Value = "  +4.17%  "

$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "43.44"
$ws.Range("E8").Value = "  +6.43%  "

$ws.Range("D9").Value = "0.332"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").Value = "0.0702"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").Value = "0.0992"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("D12").Value = "2.148.11"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.929.51"
$ws.Range("E13").Value = "  +3.96%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "11.78"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").Value = "0.684"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "4.76"
$ws.Range("E16").Value = "  +1.59%  "

$ws.Range("D17").Value = "35.507.35"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "70.95"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "242.48"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "4.80"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").Value = "2.28"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").Value = "170.95"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").Value = "2.05"
$ws.Range("E26").Value = "  +30.74%  "

$ws.Range("D27").Value = "8.24"
$ws.Range("E27").Value = "  +4.67%  "

$ws.Range("D28").Value = "17.81"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").Value = "0.0564"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("B31").Value = "BinanceUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D31").Value = "1.02"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.06"
$ws.Range("E32").Value = "  +2.19%  "

$ws.Range("D33").Value = "4.07"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("D34").Value = "0.904"
$ws.Range("E34").Value = "  +18.08%  "

$ws.Range("D35").Value = "1.75"
$ws.Range("E35").Value = "  +9.42%  "

$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  +3.89%  "

$ws.Range("D37").Value = "1.35"
$ws.Range("E37").Value = "  +10.70%  "

$ws.Range("D38").Value = "1.10"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("E39").Value = "  +3.66%  "

$ws.Range("D40").Value = "89.49"
$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("D41").Value = "1.353.65"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "15.27"
$ws.Range("E42").Value = "  +3.99%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.0592"
$ws.Range("E43").Value = "  +11.39%  "

$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "47.91"
$ws.Range("E45").Value = "  +40.84%  "

$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").Value = "12.59"
$ws.Range("E47").Value = "  +46.17%  "

$ws.Range("D48").Value = "6.70"
$ws.Range("E48").Value = "  +5.49%  "

$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("D50").Value = "2.066.44"
$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("E51").Value = "  +2.49%  "
